# Talent.xlsx: "unify the conception of DataNode, DataTable, Entity."
#
# The sheet that used to describe a generic "Property" now represents a
# "DataNode", so it is renamed accordingly. The remaining tweaks mirror the
# small layout adjustments (active cell, row heights, column widths) that
# were made in the same editing session.

$wb = $excel.ActiveWorkbook

# --- Rename the worksheet: "Property1" -> "DataNode" -----------------------
$ws = $wb.Worksheets.Item("Property1")
$ws.Name = "DataNode"

# --- Move the selection that was left active on save -----------------------
$ws.Range("D22").Select()

# --- Minor row-height touch-ups --------------------------------------------
$ws.Rows.Item(1).RowHeight = 27
$ws.Rows.Item(8).RowHeight = 27

# --- Minor column-width touch-ups -------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.16071428571429
$ws.Columns.Item(8).ColumnWidth = 25.41071428571429
